$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.627.78'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.82%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.824.46'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.55%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.008'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '1.008'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.01%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '306.99'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4645'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.66%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3604'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.27%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07136'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.17%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9018'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.57%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07767'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.79%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '19.40'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.829.02'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.80%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.262'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.50%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.328'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.13%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.40'
$ws.Range("D16").Style = "Normal"
$ws.Range("E17").Value = '  -0.05%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008562'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("E19").Value = '  +0.00%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '26.659.65'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.013'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.31%  '
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("E24").Value = '  -3.28%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '152.29'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.40%  '
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.975'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '113.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.37%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.807'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.26%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.08806'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.140'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.98%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7314'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.21%  '
$ws.Range("B33").Value = 'ARBITRUM'
$ws.Range("C33").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.140'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.42%  '
$ws.Range("B34").Value = 'Filecoin'
$ws.Range("C34").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.439'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.719'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.18%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.073'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.23%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01926'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.925'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.60%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05112'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.907'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5055'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.14%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1495'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.79%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.000'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.05%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.008'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4656'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.32%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.995'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.21%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '97.92'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.558'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.84%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05980'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '63.76'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '35.78'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.13%  '
